$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tagging")
$ws.Range("C3").Value = 3787
